# Update the "Förändrad" (changed) date column (C) for rows 2-24
# from serial date 45218 (2023-10-19) to 45221 (2023-10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
